$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R, header row 4: year 2021 — reuse Q4's style (right-aligned header, no number format)
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

# New column R, data row 5: value — reuse Q5's style (font/border/alignment), then apply the new "0.0" number format
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 102.20441221981518
$ws.Range("R5").NumberFormat = "0.0"

# Update the active selection like in the target workbook
$ws.Range("S9").Select() | Out-Null
